# Added elbow bearings into BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new row at row 8 - shifts existing rows 8+ down by one
$ws.Rows("8:8").Insert()

# 2) Populate the new row 8 with the Elbow Bearing part (C8 text/hyperlink set further below)
$ws.Range("A8").Value = "Elbow Bearing - 2788/2720 * 2"
$ws.Range("D8").Value = 44.98

# 3) The row-insert does not move the worksheet's Hyperlinks along with the cells,
#    so rebuild the hyperlink collection from scratch against the new layout.
#    (Hyperlinks collection can only be cleared as a whole in this object model.)
$ws.Hyperlinks.Delete()

$urlShoulder100 = "https://www.omc-stepperonline.com/nema-23-stepper-motor-l-56mm-gear-ratio-100-1-high-precision-planetary-gearbox-23hs22-2804s-hg100"
$urlShoulder50 = "https://www.omc-stepperonline.com/nema-23-stepper-motor-l-56mm-gear-ratio-50-1-high-precision-planetary-gearbox-23hs22-2804s-hg50"
$urlAdafruit = "https://www.amazon.com/Adafruit-Absolute-Orientation-Fusion-Breakout/dp/B017PEIGIG?crid=3J1S4PG236LH9&dib=eyJ2IjoiMSJ9.b_9piN2jpFgVdQUKUlO--RmT6WlwK-qAMknMAnS5IS0iWy-_pnTiR3lz6vGtzAZ-QU1RgFWmLARA63BbYm2iJgMJAs6M1k6BiQz6U9YyVjStz_rtSVCSfxwNcbItqNjUWQ-LD3Trq2PqxcLGdjeZ6dK9118d_v6FfAvXVBp5GxVkAJRjcc1m__B_8t4cydOF6rnBGZeh3U4aBQFniM8blFVuth9Mlcic1U3wytci0Ls.w2QFliZTY50PWL6_lk2pkuIYw4gKKLWDaw-qPCEKXlg&dib_tag=se&keywords=BNO055&qid=1748287133&sprefix=bno05%2Caps%2C183&sr=8-3"
$urlElegoo = "https://www.amazon.com/ELEGOO-Compatible-Arduino-Projects-Compliant/dp/B01H4ZLZLQ?crid=1J7XJ0X0LQPVN&dib=eyJ2IjoiMSJ9.C9li7QlUOdnawgCr8xZlUWZQDyHFpesFSB7DviVllKqGmMMk5K-VvRQhyhHA4AqyscubrPU9wtPQy7VaCUKyuO3EYxuwrXTtMddGGzmxhz2PsRzMuvE8bAxqUv28A2LO06Tke-rB36vSu9bQf09V1GC7pj5uKZzQJkLdic9dUJnzrOmyLS-h0Mapf2ito6gkx7mB0lzZBfwz7Se5A46tey_XF6O_qB5_P7jwQ0n3ykA.LeEpybFPF8KLznW1_yr7GJ-aAyT_db2aS0qKXnS1gmQ&dib_tag=se&keywords=arduino+mega&qid=1748287251&sprefix=arduino+m%2Caps%2C196&sr=8-1"
$urlBreakout = "https://www.amazon.com/Electronics-Salon-Terminal-Breakout-Arduino-MEGA-2560/dp/B07H9TJCWN?crid=3GOXBFMIMGAN5&dib=eyJ2IjoiMSJ9.3XoGINQ1IkjC7gxdNc2ljQUBQiGFRTxo6wcllu5zLJD4S5M_-a9LJcrYRsXXQWHM48vWvwebUTPYqNItuzPCHNEWmF_yUpMQKUkhu-itZOA2c2Cxf59eIFpwTYJ2f0jn0uPzD1N9kueUxKdz7t-9UC3upY-8wkvC1O3fa9MiE_LeQYs69F45wXyXzklh5Vy_W5L5hLnv2jeDA8_EWCDZzhZn5VQvzMYVrzc09GHfXxQ.CMbmQ7q0axvhf8XS3yt9uEAiVMgr_RdX5OyNJsZCLc0&dib_tag=se&keywords=arduino+mega+breakout+board&qid=1748287276&sprefix=arduino+mega+b%2Caps%2C174&sr=8-1"
$urlShoulderRotBearing = "https://www.amazon.com/dp/B0B5XQK5H1?ref=ppx_yo2ov_dt_b_fed_asin_title&th=1"
$urlShoulderRotBearing2 = "https://www.amazon.com/gp/product/B07RW366QL?smid=A1THAZDOWP300U&th=1"
$urlBelt = "https://www.mcmaster.com/6484K701/"
$urlElbowBearing = "https://www.amazon.com/uxcell-Tapered-Roller-Bearing-Width/dp/B0B5XQX2RD?crid=XEE741AKFTOD&dib=eyJ2IjoiMSJ9.DswnWR3W53Ca9cwH4GaxUPFfJg6uabw9uJwooYlFrzk7jzzlfG_MR0SsAizypRTTbxfEGlJpdYcGXDXn75oN1tfidCPGbHF7eAhr06xiHk1JEg9I9tQ5WpIR7Ey0jD_-TdYD1oo7D927GMhQ9Tolso8enodh4dCep77Z58otH--bhgN7hrswT6TSPZT1pWZTq4AFg6dcLr53SGj7WlGqL7Rqs3bs78WDBd7gsg1-6vM.f28EI9GGTBJBW6Tif3erxd-Ub8gBTqe45Gfaa4eA8E4&dib_tag=se&keywords=roller%2Btapered%2Bbearing&qid=1750119672&sprefix=roller%2Btapered%2Bbeari%2Caps%2C169&sr=8-9&th=1"

# Re-add in the same order as the original relationship ids (rId1..rId9), followed by the
# brand-new bearing link (rId10), so relationship ids line up with the target workbook.
$ws.Hyperlinks.Add($ws.Range("C2"), $urlShoulder100) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C17"), $urlAdafruit, "averageCustomerReviewsAnchor", [Type]::Missing, $urlAdafruit) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C18"), $urlElegoo, [Type]::Missing, [Type]::Missing, $urlElegoo) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C19"), $urlBreakout, [Type]::Missing, [Type]::Missing, $urlBreakout) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), $urlShoulder50) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), $urlShoulder50) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), $urlShoulderRotBearing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), $urlShoulderRotBearing2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C10"), $urlBelt) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), $urlElbowBearing, [Type]::Missing, [Type]::Missing, $urlElbowBearing) | Out-Null

# 4) Adding hyperlinks resets cell styling to a freshly-minted style; restore the shared
#    "Hyperlink" named style on every linked cell so formatting matches the original.
$linkedCells = @("C2","C3","C4","C6","C7","C8","C10","C17","C18","C19")
foreach ($addr in $linkedCells) {
    $ws.Range($addr).Style = "Hyperlink"
}
